$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.324.47'

$ws.Range('E2').Value = '  +2.25%  '

$ws.Range('D3').Value = '3.181.69'

$ws.Range('E3').Value = '  +5.36%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.69'
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').Value = '  +3.00%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '629.02'
$ws.Range('D6').Style = "Normal"

$ws.Range('E6').Value = '  +2.87%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.17'
$ws.Range('D7').Style = "Normal"

$ws.Range('E7').Value = '  +32.99%  '

$ws.Range('D10').Value = '3.178.01'

$ws.Range('E10').Value = '  +5.35%  '

$ws.Range('E11').Value = '  +15.47%  '

$ws.Range('E12').Value = '  +8.12%  '

$ws.Range('E13').Value = '  +3.24%  '

$ws.Range('E14').Value = '  +7.13%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '35.33'
$ws.Range('D15').Style = "Normal"

$ws.Range('E15').Value = '  +10.22%  '

$ws.Range('D16').Value = '90.921.69'

$ws.Range('E16').Value = '  +1.91%  '

$ws.Range('D17').Value = '3.761.19'

$ws.Range('E17').Value = '  +5.13%  '

$ws.Range('D18').Value = '3.190.19'

$ws.Range('E18').Value = '  +4.35%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.75'
$ws.Range('D19').Style = "Normal"

$ws.Range('E19').Value = '  +13.59%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.65'
$ws.Range('D20').Style = "Normal"

$ws.Range('E20').Value = '  +9.98%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '472.95'
$ws.Range('D21').Style = "Normal"

$ws.Range('E21').Value = '  +12.03%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0000212'
$ws.Range('D22').Style = "Normal"

$ws.Range('E22').Value = '  -3.74%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.19'
$ws.Range('D23').Style = "Normal"

$ws.Range('E23').Value = '  +11.78%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.32'
$ws.Range('D24').Style = "Normal"

$ws.Range('E24').Value = '  +6.06%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '96.86'
$ws.Range('D25').Style = "Normal"

$ws.Range('E25').Value = '  +18.05%  '

$ws.Range('E26').Value = '  +11.73%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.50'
$ws.Range('D27').Style = "Normal"

$ws.Range('E27').Value = '  +8.53%  '

$ws.Range('D28').Value = '3.340.57'

$ws.Range('E28').Value = '  +4.80%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = "Normal"

$ws.Range('E29').Value = '  -0.14%  '

$ws.Range('E30').Value = '  +13.13%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.164'
$ws.Range('D31').Style = "Normal"

$ws.Range('E31').Value = '  +1.26%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = "Normal"

$ws.Range('E32').Value = '  -0.02%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.96'
$ws.Range('D33').Style = "Normal"

$ws.Range('E33').Value = '  +23.50%  '

$ws.Range('B34').Value = 'Bittensor'

$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '529.76'
$ws.Range('D34').Style = "Normal"

$ws.Range('E34').Value = '  +6.27%  '

$ws.Range('B35').Value = 'Stellar'

$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.192'
$ws.Range('D35').Style = "Normal"

$ws.Range('E35').Value = '  +41.50%  '

$ws.Range('B36').Value = 'Kaspa'

$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.147'
$ws.Range('D36').Style = "Normal"

$ws.Range('E36').Value = '  +11.15%  '

$ws.Range('B37').Value = 'PancakeSwap'

$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.94'
$ws.Range('D37').Style = "Normal"

$ws.Range('E37').Value = '  +8.14%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.66'
$ws.Range('D38').Style = "Normal"

$ws.Range('E38').Value = '  -2.41%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.05'
$ws.Range('D39').Style = "Normal"

$ws.Range('E39').Value = '  +6.50%  '

$ws.Range('E40').Value = '  +6.74%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0894'
$ws.Range('D41').Style = "Normal"

$ws.Range('E41').Value = '  +28.89%  '

$ws.Range('E42').Value = '  +0.19%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.423'
$ws.Range('D43').Style = "Normal"

$ws.Range('E43').Value = '  +18.74%  '

$ws.Range('E44').Value = '  +0.03%  '

$ws.Range('E45').Value = '  +10.31%  '

$ws.Range('E46').Value = '  +0.06%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.717'
$ws.Range('D47').Style = "Normal"

$ws.Range('E47').Value = '  +22.27%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '151.51'
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').Value = '  +4.86%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.67'
$ws.Range('D49').Style = "Normal"

$ws.Range('E49').Value = '  +11.96%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '45.38'
$ws.Range('D51').Style = "Normal"
